$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 34746.207
$ws.Range("I8").Value = 250041.25
$ws.Range("J8").Value = 299
$ws.Range("K8").Value = 750123.75
$ws.Range("L8").Value = 897
$ws.Range("M8").Value = -749984.75
$ws.Range("N8").Value = -1175
$ws.Range("H17").Value = 229718.48
$ws.Range("J17").Value = 235022.1
$ws.Range("L17").Value = 705066.3
$ws.Range("N17").Value = -705402.3
$ws.Range("H18").Value = 800
$ws.Range("I18").Value = 999
$ws.Range("K18").Value = 999
$ws.Range("M18").Value = -715
$ws.Range("H32").Value = 2855
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2855
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2855
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3507
$ws.Range("H62").Value = 7426.5654
$ws.Range("I62").Value = 7448.3335
$ws.Range("K62").Value = 7448.3335
$ws.Range("M62").Value = -6824.3335
$ws.Range("H65").Value = 7426.5654
$ws.Range("I65").Value = 7448.3335
$ws.Range("K65").Value = 37241.6675
$ws.Range("M65").Value = -34121.6675
$ws.Range("H92").Value = 1349
$ws.Range("I92").Value = 1080.8422
$ws.Range("K92").Value = 1080.8422
$ws.Range("M92").Value = 167.1578
$ws.Range("H115").Value = 872.7
$ws.Range("I115").Value = 872.7
$ws.Range("K115").Value = 2618.1
$ws.Range("M115").Value = -1051.1
$ws.Range("H132").Value = 7662
$ws.Range("I132").Value = 2621.5
$ws.Range("K132").Value = 7864.5
$ws.Range("M132").Value = -5334.5
$ws.Range("H135").Value = 2189.3635
$ws.Range("I135").Value = 882
$ws.Range("J135").Value = 5675.6665
$ws.Range("K135").Value = 7938
$ws.Range("L135").Value = 51080.9985
$ws.Range("M135").Value = -5403
$ws.Range("N135").Value = -56150.9985
$ws.Range("H138").Value = 191314.81
$ws.Range("I138").Value = 53609.367
$ws.Range("J138").Value = 258402.08
$ws.Range("K138").Value = 160828.101
$ws.Range("L138").Value = 775206.24
$ws.Range("M138").Value = -155688.101
$ws.Range("N138").Value = -785486.24
$ws.Range("H141").Value = 3133.1738
$ws.Range("I141").Value = 2253.889
$ws.Range("K141").Value = 6761.667
$ws.Range("M141").Value = -1581.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6915.4854
$ws.Range("I32").Value = 6951.4126
$ws.Range("K32").Value = 6951.4126
$ws.Range("M32").Value = -6664.4126
$ws.Range("H45").Value = 2964.8
$ws.Range("I45").Value = 2556.8572
$ws.Range("K45").Value = 2556.8572
$ws.Range("M45").Value = -2179.8572
$ws.Range("H74").Value = 2942.261
$ws.Range("I74").Value = 1041.0667
$ws.Range("K74").Value = 1041.0667
$ws.Range("M74").Value = -167.0667000000001
$ws.Range("H77").Value = 2942.261
$ws.Range("I77").Value = 1041.0667
$ws.Range("K77").Value = 5205.333500000001
$ws.Range("M77").Value = -837.3335000000006
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H122").Value = 3417.6316
$ws.Range("I122").Value = 2953.5334
$ws.Range("K122").Value = 8860.600199999999
$ws.Range("M122").Value = -6410.600199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 147998
$ws.Range("J87").Value = 147998
$ws.Range("L87").Value = 147998
$ws.Range("N87").Value = -150494
$ws.Range("H90").Value = 147998
$ws.Range("J90").Value = 147998
$ws.Range("L90").Value = 443994
$ws.Range("N90").Value = -456474
$ws.Range("H105").Value = 3866.9656
$ws.Range("I105").Value = 3717
$ws.Range("J105").Value = 4260.625
$ws.Range("K105").Value = 3717
$ws.Range("L105").Value = 4260.625
$ws.Range("M105").Value = -1970
$ws.Range("N105").Value = -7754.625
$ws.Range("H107").Value = 1467.4736
$ws.Range("I107").Value = 1069.2727
$ws.Range("K107").Value = 1069.2727
$ws.Range("M107").Value = 850.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2047.7028
$ws.Range("I31").Value = 1274.4333
$ws.Range("K31").Value = 1274.4333
$ws.Range("M31").Value = -979.4332999999999
$ws.Range("H34").Value = 2047.7028
$ws.Range("I34").Value = 1274.4333
$ws.Range("K34").Value = 1274.4333
$ws.Range("M34").Value = -1072.4333
$ws.Range("H41").Value = 50000
$ws.Range("J41").Value = 50000
$ws.Range("L41").Value = 50000
$ws.Range("N41").Value = -50856
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H99").Value = 5333.1665
$ws.Range("I99").Value = 4999.75
$ws.Range("K99").Value = 4999.75
$ws.Range("M99").Value = -3501.75
$ws.Range("H103").Value = 62750
$ws.Range("I103").Value = 62750
$ws.Range("K103").Value = 62750
$ws.Range("M103").Value = -61578
$ws.Range("H126").Value = 5333.1665
$ws.Range("I126").Value = 4999.75
$ws.Range("K126").Value = 14999.25
$ws.Range("M126").Value = -12529.25
$ws.Range("H134").Value = 3036.9412
$ws.Range("I134").Value = 1437.2609
$ws.Range("J134").Value = 6381.727
$ws.Range("K134").Value = 4311.7827
$ws.Range("L134").Value = 19145.181
$ws.Range("M134").Value = -1776.7827
$ws.Range("N134").Value = -24215.181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 953.63635
$ws.Range("J5").Value = 2155.6667
$ws.Range("L5").Value = 6467.000100000001
$ws.Range("N5").Value = -6691.000100000001
$ws.Range("H12").Value = 712.23334
$ws.Range("J12").Value = 586.7222
$ws.Range("L12").Value = 1760.1666
$ws.Range("N12").Value = -2106.1666
$ws.Range("H129").Value = 67406.92999999999
$ws.Range("I129").Value = 77316.53999999999
$ws.Range("J129").Value = 2994.5
$ws.Range("K129").Value = 231949.62
$ws.Range("L129").Value = 8983.5
$ws.Range("M129").Value = -226949.62
$ws.Range("N129").Value = -18983.5
$ws.Range("H135").Value = 953.63635
$ws.Range("J135").Value = 2155.6667
$ws.Range("L135").Value = 19401.0003
$ws.Range("N135").Value = -24471.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4345.926
$ws.Range("I80").Value = 4412.125
$ws.Range("J80").Value = 4249.636
$ws.Range("K80").Value = 4412.125
$ws.Range("L80").Value = 4249.636
$ws.Range("M80").Value = -3414.125
$ws.Range("N80").Value = -6245.636
$ws.Range("H83").Value = 4345.926
$ws.Range("I83").Value = 4412.125
$ws.Range("J83").Value = 4249.636
$ws.Range("K83").Value = 22060.625
$ws.Range("L83").Value = 21248.18
$ws.Range("M83").Value = -17068.625
$ws.Range("N83").Value = -31232.18
$ws.Range("H97").Value = 973.57574
$ws.Range("I97").Value = 505.31818
$ws.Range("J97").Value = 1910.091
$ws.Range("K97").Value = 505.31818
$ws.Range("L97").Value = 1910.091
$ws.Range("M97").Value = -9.318179999999984
$ws.Range("N97").Value = -2902.091
$ws.Range("H102").Value = 20643.852
$ws.Range("I102").Value = 1784.05
$ws.Range("J102").Value = 74529
$ws.Range("K102").Value = 1784.05
$ws.Range("L102").Value = 74529
$ws.Range("M102").Value = -162.05
$ws.Range("N102").Value = -77773
$ws.Range("H107").Value = 770.8889
$ws.Range("I107").Value = 654.875
$ws.Range("J107").Value = 863.7
$ws.Range("K107").Value = 654.875
$ws.Range("L107").Value = 863.7
$ws.Range("M107").Value = 1265.125
$ws.Range("N107").Value = -4703.7
$ws.Range("H113").Value = 1557.7142
$ws.Range("I113").Value = 1571.7693
$ws.Range("K113").Value = 1571.7693
$ws.Range("M113").Value = 598.2307000000001
$ws.Range("H132").Value = 5131278
$ws.Range("I132").Value = 6175781
$ws.Range("J132").Value = 3717.182
$ws.Range("K132").Value = 18527343
$ws.Range("L132").Value = 11151.546
$ws.Range("M132").Value = -18524813
$ws.Range("N132").Value = -16211.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5615.5864
$ws.Range("I40").Value = 5598.9546
$ws.Range("K40").Value = 5598.9546
$ws.Range("M40").Value = -5462.9546
$ws.Range("H122").Value = 3792.9285
$ws.Range("I122").Value = 2711.4443
$ws.Range("J122").Value = 5739.6
$ws.Range("K122").Value = 8134.3329
$ws.Range("L122").Value = 17218.8
$ws.Range("M122").Value = -5684.3329
$ws.Range("N122").Value = -22118.8
$ws.Range("H123").Value = 83476.336
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 83476.336
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 83476.336
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -93276.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1572.4
$ws.Range("I23").Value = 740.5
$ws.Range("K23").Value = 740.5
$ws.Range("M23").Value = -511.5
$ws.Range("H43").Value = 23290
$ws.Range("I43").Value = 9999
$ws.Range("K43").Value = 9999
$ws.Range("M43").Value = -9850
$ws.Range("H96").Value = 3424
$ws.Range("J96").Value = 3708.8
$ws.Range("L96").Value = 3708.8
$ws.Range("N96").Value = -6454.8
$ws.Range("H122").Value = 2174.074
$ws.Range("I122").Value = 1989.0526
$ws.Range("J122").Value = 2613.5
$ws.Range("K122").Value = 5967.1578
$ws.Range("L122").Value = 7840.5
$ws.Range("M122").Value = -3517.1578
$ws.Range("N122").Value = -12740.5
$ws.Range("H132").Value = 2331.423
$ws.Range("I132").Value = 2302.8948
$ws.Range("K132").Value = 6908.6844
$ws.Range("M132").Value = -4378.6844
$ws.Range("H136").Value = 9001.333
$ws.Range("I136").Value = 10812.096
$ws.Range("K136").Value = 32436.288
$ws.Range("M136").Value = -29886.288
